$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,13
$data[0,0] = 1.02
$data[0,1] = 1.051534985894685
$data[0,2] = 1.060203582676322
$data[0,3] = 1.059213121012961
$data[0,4] = 1.071713674048925
$data[0,5] = 1
$data[0,6] = $null
$data[0,7] = 1.052310826394575
$data[0,8] = 1.056562347819889
$data[0,9] = 1.062931000765941
$data[0,10] = 1.061943238859832
$data[0,11] = 1.074410114180239
$data[0,12] = 1.058062787543189
$data[1,0] = 1.02
$data[1,1] = 1.052424414698844
$data[1,2] = 1.060924545348348
$data[1,3] = 1.059994967393822
$data[1,4] = 1.072553799750207
$data[1,5] = 1
$data[1,6] = $null
$data[1,7] = 1.052579425035528
$data[1,8] = 1.057102191549707
$data[1,9] = 1.063466578983021
$data[1,10] = 1.062539352416474
$data[1,11] = 1.075066787708864
$data[1,12] = 1.058603397913023
$data[2,0] = 1.02
$data[2,1] = 1.05300053871594
$data[2,2] = 1.061391557297792
$data[2,3] = 1.060501752128408
$data[2,4] = 1.073098315647763
$data[2,5] = 1
$data[2,6] = $null
$data[2,7] = 1.052752257466571
$data[2,8] = 1.05745143974822
$data[2,9] = 1.063812960939062
$data[2,10] = 1.062925293679281
$data[2,11] = 1.075491946127797
$data[2,12] = 1.058953142084038
$data[3,0] = 1.02
$data[3,1] = 1.053242884411845
$data[3,2] = 1.061588007360402
$data[3,3] = 1.060715012872876
$data[3,4] = 1.073327443403826
$data[3,5] = 1
$data[3,6] = $null
$data[3,7] = 1.05282468349238
$data[3,8] = 1.05759824648711
$data[3,9] = 1.063958537227258
$data[3,10] = 1.063087593729799
$data[3,11] = 1.075670740304641
$data[3,12] = 1.059100157305336
$data[4,0] = 1.02
$data[4,1] = 1.053283583672968
$data[4,2] = 1.061620999075215
$data[4,3] = 1.06075083246351
$data[4,4] = 1.073365927416251
$data[4,5] = 1
$data[4,6] = $null
$data[4,7] = 1.052836830470775
$data[4,8] = 1.057622894940459
$data[4,9] = 1.06398297758633
$data[4,10] = 1.063114847535354
$data[4,11] = 1.075700763966689
$data[4,12] = 1.059124840762316
$data[5,0] = 1.02
$data[5,1] = 1.053003776390299
$data[5,2] = 1.061394181810647
$data[5,3] = 1.060504600912033
$data[5,4] = 1.073101376428302
$data[5,5] = 1
$data[5,6] = $null
$data[5,7] = 1.052753226142111
$data[5,8] = 1.057453401456168
$data[5,9] = 1.063814906305878
$data[5,10] = 1.062927462144377
$data[5,11] = 1.075494334959968
$data[5,12] = 1.058955106577836
$data[6,0] = 1.02
$data[6,1] = 1.051835446557093
$data[6,2] = 1.060447130666946
$data[6,3] = 1.059477166862653
$data[6,4] = 1.071997411360675
$data[6,5] = 1
$data[6,6] = $null
$data[6,7] = 1.052401800698094
$data[6,8] = 1.056744803636049
$data[6,9] = 1.06311203697644
$data[6,10] = 1.062144652705278
$data[6,11] = 1.07463198808875
$data[6,12] = 1.058245502467536
$data[7,0] = 1.02
$data[7,1] = 1.049781391704661
$data[7,2] = 1.058782215819731
$data[7,3] = 1.05767349589791
$data[7,4] = 1.070059044221221
$data[7,5] = 1
$data[7,6] = $null
$data[7,7] = 1.051775158615852
$data[7,8] = 1.055495705841967
$data[7,9] = 1.061872223241136
$data[7,10] = 1.060766955940122
$data[7,11] = 1.073114378362092
$data[7,12] = 1.056994630811312
$data[8,0] = 1.02
$data[8,1] = 1.048415263091769
$data[8,2] = 1.057675001332809
$data[8,3] = 1.056475722088942
$data[8,4] = 1.068771586610706
$data[8,5] = 1
$data[8,6] = $null
$data[8,7] = 1.051352478213645
$data[8,8] = 1.054662738241775
$data[8,9] = 1.061044904122311
$data[8,10] = 1.059849723830014
$data[8,11] = 1.072104042148398
$data[8,12] = 1.056160480301586
$data[9,0] = 1.02
$data[9,1] = 1.04782450127505
$data[9,2] = 1.057196233449806
$data[9,3] = 1.055958202822289
$data[9,4] = 1.068215261314317
$data[9,5] = 1
$data[9,6] = $null
$data[9,7] = 1.051168295972013
$data[9,8] = 1.054302013300096
$data[9,9] = 1.06068649757222
$data[9,10] = 1.059452861404173
$data[9,11] = 1.071666906993197
$data[9,12] = 1.055799243089113
$data[10,0] = 1.02
$data[10,1] = 1.04760518479406
$data[10,2] = 1.057018499083271
$data[10,3] = 1.05576614383466
$data[10,4] = 1.068008792279148
$data[10,5] = 1
$data[10,6] = $null
$data[10,7] = 1.051099709148707
$data[10,8] = 1.054168018520321
$data[10,9] = 1.060553344866993
$data[10,10] = 1.059305496225505
$data[10,11] = 1.071504589176187
$data[10,12] = 1.055665058021382
$data[11,0] = 1.02
$data[11,1] = 1.047652223526322
$data[11,2] = 1.057056619095152
$data[11,3] = 1.055807333386563
$data[11,4] = 1.068053072641922
$data[11,5] = 1
$data[11,6] = $null
$data[11,7] = 1.051114429086746
$data[11,8] = 1.054196761083817
$data[11,9] = 1.060581907664949
$data[11,10] = 1.059337104399081
$data[11,11] = 1.071539404447445
$data[11,12] = 1.055693841402615
$data[12,0] = 1.02
$data[12,1] = 1.047806370075264
$data[12,2] = 1.057181539794502
$data[12,3] = 1.055942323673969
$data[12,4] = 1.068198190939984
$data[12,5] = 1
$data[12,6] = $null
$data[12,7] = 1.051162630095418
$data[12,8] = 1.05429093735451
$data[12,9] = 1.060675491621606
$data[12,10] = 1.059440679181905
$data[12,11] = 1.07165348864885
$data[12,12] = 1.055788151414414
$data[13,0] = 1.02
$data[13,1] = 1.047901360631792
$data[13,2] = 1.057258521039646
$data[13,3] = 1.05602551831805
$data[13,4] = 1.068287626350083
$data[13,5] = 1
$data[13,6] = $null
$data[13,7] = 1.051192305381626
$data[13,8] = 1.054348961773693
$data[13,9] = 1.060733148570363
$data[13,10] = 1.059504501318522
$data[13,11] = 1.071723786843584
$data[13,12] = 1.055846258234928
$data[14,0] = 1.02
$data[14,1] = 1.048454486583298
$data[14,2] = 1.057706789699547
$data[14,3] = 1.056510091982661
$data[14,4] = 1.06880853252265
$data[14,5] = 1
$data[14,6] = $null
$data[14,7] = 1.051364677432225
$data[14,8] = 1.054686677534081
$data[14,9] = 1.061068686871338
$data[14,10] = 1.059876068826513
$data[14,11] = 1.072133060809533
$data[14,12] = 1.056184453590432
$data[15,0] = 1.02
$data[15,1] = 1.048801657961614
$data[15,2] = 1.057988155563027
$data[15,3] = 1.056814354580277
$data[15,4] = 1.069135593005298
$data[15,5] = 1
$data[15,6] = $null
$data[15,7] = 1.051472492150439
$data[15,8] = 1.054898506652675
$data[15,9] = 1.061279116055618
$data[15,10] = 1.06010922595555
$data[15,11] = 1.072389881538197
$data[15,12] = 1.056396583530672
$data[16,0] = 1.02
$data[16,1] = 1.049004232400331
$data[16,2] = 1.058152335396074
$data[16,3] = 1.056991934206188
$data[16,4] = 1.069326472913695
$data[16,5] = 1
$data[16,6] = $null
$data[16,7] = 1.051535266789021
$data[16,8] = 1.055022058621609
$data[16,9] = 1.061401839166154
$data[16,10] = 1.060245251914059
$data[16,11] = 1.072539714113186
$data[16,12] = 1.056520310957573
$data[17,0] = 1.02
$data[17,1] = 1.049073317768404
$data[17,2] = 1.058208327260983
$data[17,3] = 1.057052502603276
$data[17,4] = 1.069391576812292
$data[17,5] = 1
$data[17,6] = $null
$data[17,7] = 1.051556652321503
$data[17,8] = 1.055064185871875
$data[17,9] = 1.061443681708728
$data[17,10] = 1.060291638187737
$data[17,11] = 1.072590808767077
$data[17,12] = 1.056562498033367
$data[18,0] = 1.02
$data[18,1] = 1.048764401964708
$data[18,2] = 1.057957961048551
$data[18,3] = 1.056781698858158
$data[18,4] = 1.069100491009923
$data[18,5] = 1
$data[18,6] = $null
$data[18,7] = 1.051460936216467
$data[18,8] = 1.054875779836493
$data[18,9] = 1.061256540710946
$data[18,10] = 1.060084207356613
$data[18,11] = 1.072362323640153
$data[18,12] = 1.056373824439804
$data[19,0] = 1.02
$data[19,1] = 1.047760974468605
$data[19,2] = 1.057144750950999
$data[19,3] = 1.055902567672474
$data[19,4] = 1.068155452374054
$data[19,5] = 1
$data[19,6] = $null
$data[19,7] = 1.051148440878571
$data[19,8] = 1.054263204931517
$data[19,9] = 1.060647934146481
$data[19,10] = 1.059410177676162
$data[19,11] = 1.071619892202178
$data[19,12] = 1.0557603796082
$data[20,0] = 1.02
$data[20,1] = 1.047130767609747
$data[20,2] = 1.056634040657914
$data[20,3] = 1.055350811541701
$data[20,4] = 1.06756228218625
$data[20,5] = 1
$data[20,6] = $null
$data[20,7] = 1.050950960156365
$data[20,8] = 1.053878023696302
$data[20,9] = 1.0602651376103
$data[20,10] = 1.058986662199462
$data[20,11] = 1.071153407193493
$data[20,12] = 1.055374651371449
$data[21,0] = 1.02
$data[21,1] = 1.047464786508215
$data[21,2] = 1.056904721626175
$data[21,3] = 1.055643213624145
$data[21,4] = 1.067876636285862
$data[21,5] = 1
$data[21,6] = $null
$data[21,7] = 1.051055743193262
$data[21,8] = 1.054082218099728
$data[21,9] = 1.060468078277703
$data[21,10] = 1.059211149367581
$data[21,11] = 1.071400669733991
$data[21,12] = 1.055579135754351
$data[22,0] = 1.02
$data[22,1] = 1.048781236113697
$data[22,2] = 1.057971604452746
$data[22,3] = 1.056796454238608
$data[22,4] = 1.069116351747744
$data[22,5] = 1
$data[22,6] = $null
$data[22,7] = 1.051466158191329
$data[22,8] = 1.054886049120109
$data[22,9] = 1.061266741589121
$data[22,10] = 1.060095512094033
$data[22,11] = 1.072374775765074
$data[22,12] = 1.05638410830698
$data[23,0] = 1.02
$data[23,1] = 1.050311849199952
$data[23,2] = 1.059212162552886
$data[23,3] = 1.058138971818036
$data[23,4] = 1.070559322684477
$data[23,5] = 1
$data[23,6] = $null
$data[23,7] = 1.051938031045802
$data[23,8] = 1.055818674397763
$data[23,9] = 1.062192887264968
$data[23,10] = 1.061122912242699
$data[23,11] = 1.073506475957713
$data[23,12] = 1.057318058019503

$ws.Range("B2:N25").Value = $data
